$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 62.85906733333334
$ws.Range("H2").Value = 188.577202
$ws.Range("I2").Value = 0.145580545806332
$ws.Range("J2").Value = 0.145580545806332
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.751166666666666
$ws.Range("N2").Value = 17.2535
$ws.Range("O2").Value = 0.7405222614421495
$ws.Range("P2").Value = 0.7405222614421495
$ws.Range("Q2").Value = 361.5129727452222
$ws.Range("R2").Value = 3253.616754707
$ws.Range("S2").Value = 0.1078056350024874
$ws.Range("T2").Value = 0.1078056350024874

# Row 3
$ws.Range("G3").Value = 62.85906733333334
$ws.Range("H3").Value = 188.577202
$ws.Range("I3").Value = 0.145580545806332
$ws.Range("J3").Value = 0.145580545806332
$ws.Range("O3").Value = 0.07337387367415998
$ws.Range("P3").Value = 0.07337387367416
$ws.Range("Q3").Value = 35.82013475478777
$ws.Range("R3").Value = 322.38121279309
$ws.Range("S3").Value = 0.01068180857740906
$ws.Range("T3").Value = 0.01068180857740906

# Row 4
$ws.Range("G4").Value = 62.85906733333334
$ws.Range("H4").Value = 188.577202
$ws.Range("I4").Value = 0.145580545806332
$ws.Range("J4").Value = 0.145580545806332
$ws.Range("M4").Value = 1.445350666666667
$ws.Range("N4").Value = 4.336052
$ws.Range("O4").Value = 0.1861038648836906
$ws.Range("P4").Value = 0.1861038648836906
$ws.Range("Q4").Value = 90.85339487627824
$ws.Range("R4").Value = 817.6805538865041
$ws.Range("S4").Value = 0.02709310222643553
$ws.Range("T4").Value = 0.02709310222643553

# Row 5
$ws.Range("I5").Value = 0.331880415407135
$ws.Range("J5").Value = 0.331880415407135
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.751166666666666
$ws.Range("N5").Value = 17.2535
$ws.Range("O5").Value = 0.7405222614421495
$ws.Range("P5").Value = 0.7405222614421495
$ws.Range("Q5").Value = 824.1422293426667
$ws.Range("R5").Value = 7417.280064084001
$ws.Range("S5").Value = 0.2457648357456516
$ws.Range("T5").Value = 0.2457648357456516

# Row 6
$ws.Range("I6").Value = 0.331880415407135
$ws.Range("J6").Value = 0.331880415407135
$ws.Range("O6").Value = 0.07337387367415998
$ws.Range("P6").Value = 0.07337387367416
$ws.Range("S6").Value = 0.02435135167501086
$ws.Range("T6").Value = 0.02435135167501087

# Row 7
$ws.Range("I7").Value = 0.331880415407135
$ws.Range("J7").Value = 0.331880415407135
$ws.Range("M7").Value = 1.445350666666667
$ws.Range("N7").Value = 4.336052
$ws.Range("O7").Value = 0.1861038648836906
$ws.Range("P7").Value = 0.1861038648836906
$ws.Range("Q7").Value = 207.1187620961387
$ws.Range("R7").Value = 1864.068858865248
$ws.Range("S7").Value = 0.06176422798647255
$ws.Range("T7").Value = 0.06176422798647255

# Row 8
$ws.Range("G8").Value = 144.7357836666667
$ws.Range("H8").Value = 434.207351
$ws.Range("I8").Value = 0.3352056477733801
$ws.Range("J8").Value = 0.3352056477733801
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.751166666666666
$ws.Range("N8").Value = 17.2535
$ws.Range("O8").Value = 0.7405222614421495
$ws.Range("P8").Value = 0.7405222614421495
$ws.Range("Q8").Value = 832.3996144976111
$ws.Range("R8").Value = 7491.5965304785
$ws.Range("S8").Value = 0.248227244337324
$ws.Range("T8").Value = 0.2482272443373241

# Row 9
$ws.Range("G9").Value = 144.7357836666667
$ws.Range("H9").Value = 434.207351
$ws.Range("I9").Value = 0.3352056477733801
$ws.Range("J9").Value = 0.3352056477733801
$ws.Range("O9").Value = 0.07337387367415998
$ws.Range("P9").Value = 0.07337387367416
$ws.Range("Q9").Value = 82.47744509614388
$ws.Range("R9").Value = 742.297005865295
$ws.Range("S9").Value = 0.02459533685458896
$ws.Range("T9").Value = 0.02459533685458896

# Row 10
$ws.Range("G10").Value = 144.7357836666667
$ws.Range("H10").Value = 434.207351
$ws.Range("I10").Value = 0.3352056477733801
$ws.Range("J10").Value = 0.3352056477733801
$ws.Range("M10").Value = 1.445350666666667
$ws.Range("N10").Value = 4.336052
$ws.Range("O10").Value = 0.1861038648836906
$ws.Range("P10").Value = 0.1861038648836906
$ws.Range("Q10").Value = 209.1939614131391
$ws.Range("R10").Value = 1882.745652718252
$ws.Range("S10").Value = 0.06238306658146711
$ws.Range("T10").Value = 0.06238306658146712

# Row 11
$ws.Range("G11").Value = 80.88719666666667
$ws.Range("H11").Value = 242.66159
$ws.Range("I11").Value = 0.187333391013153
$ws.Range("J11").Value = 0.187333391013153
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.751166666666666
$ws.Range("N11").Value = 17.2535
$ws.Range("O11").Value = 0.7405222614421495
$ws.Range("P11").Value = 0.7405222614421495
$ws.Range("Q11").Value = 465.1957492294444
$ws.Range("R11").Value = 4186.761743065
$ws.Range("S11").Value = 0.1387245463566865
$ws.Range("T11").Value = 0.1387245463566865

# Row 12
$ws.Range("G12").Value = 80.88719666666667
$ws.Range("H12").Value = 242.66159
$ws.Range("I12").Value = 0.187333391013153
$ws.Range("J12").Value = 0.187333391013153
$ws.Range("O12").Value = 0.07337387367415998
$ws.Range("P12").Value = 0.07337387367416
$ws.Range("Q12").Value = 46.09343420850556
$ws.Range("R12").Value = 414.84090787655
$ws.Range("S12").Value = 0.01374537656715111
$ws.Range("T12").Value = 0.01374537656715111

# Row 13
$ws.Range("G13").Value = 80.88719666666667
$ws.Range("H13").Value = 242.66159
$ws.Range("I13").Value = 0.187333391013153
$ws.Range("J13").Value = 0.187333391013153
$ws.Range("M13").Value = 1.445350666666667
$ws.Range("N13").Value = 4.336052
$ws.Range("O13").Value = 0.1861038648836906
$ws.Range("P13").Value = 0.1861038648836906
$ws.Range("Q13").Value = 116.9103636269645
$ws.Range("R13").Value = 1052.19327264268
$ws.Range("S13").Value = 0.03486346808931541
$ws.Range("T13").Value = 0.03486346808931541
